# Reorder the "Recorded By" (column G) attendee lists for the listed rows.
# Each cell keeps the same set of names, just in a different order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "Dr. Mohammad El-Tanany, Dr. Rana Abo-Zaid, Dr. Servinaz Sayed Mohammad, Dr. Nesma, Dr. Nahla Nagiub"
$ws.Range("G3").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Rana Abo-Zaid, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Asmaa Reda, Dr. Veronia Rafat"
$ws.Range("G4").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Rana Abo-Zaid, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Gehan Adel, Dr. Majorelle Magdy, Dr. Manar Montaser"
$ws.Range("G5").Value = "Dr. Lamiaa Ossama, Dr. Menna tu'Alllah Mohammad, Dr. Nada Gouda, Dr. Fatma Elhady, Dr. Amera Ahmad Saad"
$ws.Range("G6").Value = "Dr. Lamiaa Ossama, Dr. Menna tu'Alllah Mohammad, Dr. Fatma Elhady, Dr. Kerelos Zareef, Dr. Abeer Ragab, Dr. Amera Ahmad Saad, Dr. Nada Mohammad"
$ws.Range("G7").Value = "Dr. Arwa Elnagar, Dr. Shimaa Ashraf, Dr. Aya Saeed"
$ws.Range("G8").Value = "Dr. Marwa Mustafa, Dr. Madeha Saeed, Dr. Amira Ibrahim, Dr. Dalia Mohammad Abd Al-Salam, Dr. Dina Adel"
$ws.Range("G9").Value = "Dr. Arwa Al-Sayed, Dr. Merna Said, Dr. Yasmeena Fattoh, Dr. Madeha Saeed, Dr. Maryam Ahmad"
$ws.Range("G10").Value = "Dr. Al-Shimaa Khaled, Dr. Mohammad Safwat"
$ws.Range("G11").Value = "D Mariam E. Mohammad, Dr. Mona Ibrahim Hussein, Dr. Heba Al-Sayed Mohammad"
$ws.Range("G12").Value = "D Mariam E. Mohammad, Dr. Mona Ibrahim Hussein, Dr. Heba Al-Sayed Mohammad"
$ws.Range("G13").Value = "Dr. Marian Samir, Dr. Aya Alaa-Eldein, Dr. Manarst Al-Eslam"
$ws.Range("G15").Value = "Dr. Amr Saeed, Dr. Walaa Ghanima"
$ws.Range("G16").Value = "Dr. Nardine, Dr. Wafaa Ebida, Dr. Neveen Nashaat, Dr. Aya Hanafy, Dr. Eman Samir Gabry, Dr. Salma Hassan, Dr. Remon, Dr. Abdullah El-Agrody"
$ws.Range("G17").Value = "Dr. Wafaa Ebida, Dr. Yasmin, Dr. Neveen Nashaat, Dr. Marina Sorial, Dr. Eman Samir Gabry"
$ws.Range("G18").Value = "Dr. Mohammad El-Tanany, Dr. Rana Abo-Zaid, Dr. Servinaz Sayed Mohammad, Dr. Nesma, Dr. Nahla Nagiub"
$ws.Range("G19").Value = "Dr. Mohammad El-Tanany, Dr. Amira Sobhy, Dr. Servinaz Sayed Mohammad, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Veronia Rafat"
$ws.Range("G20").Value = "Dr. Nourhan Mahmoud, Dr. Mohammad El-Tanany, Dr. Amira Sobhy, Dr. Heba Mahmoud Ali, Dr. Nesma, Dr. Servinaz Sayed Mohammad, Dr. Asmaa Reda"
$ws.Range("G21").Value = "Dr. Lamiaa Ossama, Dr. Menna tu'Alllah Mohammad, Dr. Nada Gouda, Dr. Fatma Elhady, Dr. Amera Ahmad Saad"
$ws.Range("G22").Value = "Dr. Lamiaa Ossama, Dr. Menna tu'Alllah Mohammad, Dr. Fatma Elhady, Dr. Kerelos Zareef, Dr. Abeer Ragab, Dr. Amera Ahmad Saad, Dr. Nada Mohammad"
$ws.Range("G23").Value = "Dr. Arwa Elnagar, Dr. Shimaa Ashraf, Dr. Aya Saeed"
$ws.Range("G24").Value = "Dr. Marwa Mustafa, Dr. Madeha Saeed, Dr. Amira Ibrahim, Dr. Dalia Mohammad Abd Al-Salam, Dr. Dina Adel"
$ws.Range("G25").Value = "Dr. Arwa Al-Sayed, Dr. Merna Said, Dr. Yasmeena Fattoh, Dr. Madeha Saeed, Dr. Maryam Ahmad"
$ws.Range("G26").Value = "Dr. Al-Shimaa Khaled, Dr. Mohammad Safwat"
$ws.Range("G27").Value = "D Mariam E. Mohammad, Dr. Mona Ibrahim Hussein, Dr. Heba Al-Sayed Mohammad"
$ws.Range("G28").Value = "D Mariam E. Mohammad, Dr. Mona Ibrahim Hussein, Dr. Heba Al-Sayed Mohammad"
$ws.Range("G29").Value = "Dr. Marian Samir, Dr. Aya Alaa-Eldein, Dr. Manarst Al-Eslam"
$ws.Range("G31").Value = "Dr. Amr Saeed, Dr. Walaa Ghanima"
$ws.Range("G32").Value = "Dr. Nardine, Dr. Wafaa Ebida, Dr. Neveen Nashaat, Dr. Aya Hanafy, Dr. Eman Samir Gabry, Dr. Salma Hassan, Dr. Remon, Dr. Abdullah El-Agrody"
$ws.Range("G33").Value = "Dr. Wafaa Ebida, Dr. Yasmin, Dr. Neveen Nashaat, Dr. Marina Sorial, Dr. Eman Samir Gabry"
$ws.Range("G34").Value = "Dr. Menna tuâ€™Allah Medhat, Administrator, Dr. Rana Abo-Zaid, Dr. Amira Sobhy, Dr. Servinaz Sayed Mohammad, Dr. Asmaa Reda, Dr. Gehan Adel, Dr. Nahla Nagiub, Dr. Veronia Rafat"
$ws.Range("G35").Value = "Administrator, Dr. Mohammad El-Tanany, Dr. Rana Abo-Zaid, Dr. Heba Mahmoud Ali, Dr. Servinaz Sayed Mohammad, Dr. Eman Tantawi, Dr. Gehan Adel"
$ws.Range("G36").Value = "Dr. Shimaa Ahmad Mekki, Dr. Rana Abo-Zaid, Dr. Amira Sobhy, Dr. Alshimaa Atef, Dr. Heba Mahmoud Ali, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat"
$ws.Range("G37").Value = "Dr. Lamiaa Ossama, Dr. Nada Gouda, Dr. Fatma Elhady, Dr. Kerelos Zareef, Dr. Abeer Ragab, Dr. Nada Mohammad"
$ws.Range("G40").Value = "Dr. Eman M. Abo-Sakaya, Dr. Merna Said, Dr. Sara Atawia, Dr. Basma Hamed, Dr. Nourhan Osama, Dr. Amany Raafat, Dr. Yasmeena Fattoh, Dr. Madeha Saeed, Dr. Merna Mahrous, Dr. Marina Youhanna, Dr. Maryam Ahmad, Dr. Nahed Mosaad, Dr. Mai Mustafa"
$ws.Range("G41").Value = "Dr. Eman M. Abo-Sakaya, Dr. Esraa Mostafa, Dr. Nourhan Osama, Dr. Amany Raafat, Dr. Sarah Abdelmohsen, Dr. Merna Mahrous, Dr. Amira Ibrahim, Dr. Nadia Mostafa, Dr. Dina Adel, Dr. Maryam Ahmad"
$ws.Range("G47").Value = "Dr. Afaf Abdallah, Dr. Aya Alaa-Eldein"
$ws.Range("G48").Value = "Dr. Maryam Ashraf, Dr. Remon"
$ws.Range("G49").Value = "Dr. Ola Abd Al-Fattah, Dr. Naema Gomaa, Dr. Neveen Nashaat, Dr. Monica, Dr. Eman Samir Gabry, Dr. Remon"
$ws.Range("G50").Value = "Dr. Menna tuâ€™Allah Medhat, Administrator, Dr. Rana Abo-Zaid, Dr. Amira Sobhy, Dr. Servinaz Sayed Mohammad, Dr. Asmaa Reda, Dr. Gehan Adel, Dr. Nahla Nagiub, Dr. Veronia Rafat"
$ws.Range("G51").Value = "Administrator, Dr. Mohammad El-Tanany, Dr. Rana Abo-Zaid, Dr. Heba Mahmoud Ali, Dr. Servinaz Sayed Mohammad, Dr. Eman Tantawi, Dr. Gehan Adel"
$ws.Range("G52").Value = "Dr. Shimaa Ahmad Mekki, Dr. Rana Abo-Zaid, Dr. Amira Sobhy, Dr. Alshimaa Atef, Dr. Heba Mahmoud Ali, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat"
$ws.Range("G53").Value = "Dr. Lamiaa Ossama, Dr. Nada Gouda, Dr. Fatma Elhady, Dr. Kerelos Zareef, Dr. Abeer Ragab, Dr. Nada Mohammad"
$ws.Range("G56").Value = "Dr. Eman M. Abo-Sakaya, Dr. Merna Said, Dr. Sara Atawia, Dr. Basma Hamed, Dr. Nourhan Osama, Dr. Amany Raafat, Dr. Yasmeena Fattoh, Dr. Madeha Saeed, Dr. Merna Mahrous, Dr. Marina Youhanna, Dr. Maryam Ahmad, Dr. Nahed Mosaad, Dr. Mai Mustafa"
$ws.Range("G57").Value = "Dr. Eman M. Abo-Sakaya, Dr. Esraa Mostafa, Dr. Nourhan Osama, Dr. Amany Raafat, Dr. Sarah Abdelmohsen, Dr. Merna Mahrous, Dr. Amira Ibrahim, Dr. Nadia Mostafa, Dr. Dina Adel, Dr. Maryam Ahmad"
$ws.Range("G63").Value = "Dr. Afaf Abdallah, Dr. Aya Alaa-Eldein"
$ws.Range("G64").Value = "Dr. Maryam Ashraf, Dr. Remon"
$ws.Range("G65").Value = "Dr. Ola Abd Al-Fattah, Dr. Naema Gomaa, Dr. Neveen Nashaat, Dr. Monica, Dr. Eman Samir Gabry, Dr. Remon"
$ws.Range("G66").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Nourhan Mahmoud, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Gehan Adel, Dr. Nahla Nagiub, Dr. Veronia Rafat"
$ws.Range("G67").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Rana Abo-Zaid, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Asmaa Reda, Dr. Veronia Rafat"
$ws.Range("G68").Value = "Dr. Nourhan Mahmoud, Dr. Amira Sobhy, Dr. Alshimaa Atef, Dr. Eman Tantawi, Dr. Veronia Rafat"
$ws.Range("G69").Value = "Dr. Lamiaa Ossama, Dr. Menna tu'Alllah Mohammad, Dr. Abeer Ragab, Dr. Amera Ahmad Saad, Dr. Nada Mohammad"
$ws.Range("G70").Value = "Dr. Fatma Elhady, Dr. Amera Ahmad Saad, Dr. Nada Gouda"
$ws.Range("G71").Value = "Dr. Sara Nabil, Dr. Nourhan Mohammad"
$ws.Range("G72").Value = "Dr. Eman M. Abo-Sakaya, Dr. Merna Said, Dr. Sara Atawia, Dr. Basma Hamed, Dr. Nourhan Osama, Dr. Amany Raafat, Dr. Yasmeena Fattoh, Dr. Madeha Saeed, Dr. Merna Mahrous, Dr. Marina Youhanna, Dr. Maryam Ahmad, Dr. Nahed Mosaad, Dr. Mai Mustafa"
$ws.Range("G73").Value = "Dr. Arwa Al-Sayed, Dr. Merna Said, Dr. Esraa Mostafa, Dr. Madeha Saeed, Dr. Dalia Mohammad Abd Al-Salam, Dr. Nahed Mosaad"
$ws.Range("G74").Value = "Dr. Al-Shimaa Khaled, Dr. Mohammad Safwat"
$ws.Range("G79").Value = "Dr. Amr Saeed, Dr. Walaa Ghanima"
$ws.Range("G80").Value = "Dr. Ola Abd Al-Fattah, Dr. Eman Mohammad Al, Dr. Neveen Nashaat, Dr. Aya Hanafy, Dr. Salma Hassan, Dr. Marina Atef"
$ws.Range("G81").Value = "Dr. Wafaa Ebida, Dr. Yasmin, Dr. Neveen Nashaat, Dr. Marina Sorial, Dr. Eman Samir Gabry"
$ws.Range("G82").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Nourhan Mahmoud, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Gehan Adel, Dr. Nahla Nagiub, Dr. Veronia Rafat"
$ws.Range("G83").Value = "Dr. Mohammad El-Tanany, Dr. Amira Sobhy, Dr. Servinaz Sayed Mohammad, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Veronia Rafat"
$ws.Range("G84").Value = "Dr. Nourhan Mahmoud, Dr. Mohammad El-Tanany, Dr. Amira Sobhy, Dr. Heba Mahmoud Ali, Dr. Nesma, Dr. Servinaz Sayed Mohammad, Dr. Asmaa Reda"
$ws.Range("G85").Value = "Dr. Lamiaa Ossama, Dr. Menna tu'Alllah Mohammad, Dr. Abeer Ragab, Dr. Amera Ahmad Saad, Dr. Nada Mohammad"
$ws.Range("G86").Value = "Dr. Fatma Elhady, Dr. Amera Ahmad Saad, Dr. Nada Gouda"
$ws.Range("G87").Value = "Dr. Sara Nabil, Dr. Nourhan Mohammad"
$ws.Range("G88").Value = "Dr. Eman M. Abo-Sakaya, Dr. Merna Said, Dr. Sara Atawia, Dr. Basma Hamed, Dr. Nourhan Osama, Dr. Amany Raafat, Dr. Yasmeena Fattoh, Dr. Madeha Saeed, Dr. Merna Mahrous, Dr. Marina Youhanna, Dr. Maryam Ahmad, Dr. Nahed Mosaad, Dr. Mai Mustafa"
$ws.Range("G89").Value = "Dr. Arwa Al-Sayed, Dr. Merna Said, Dr. Esraa Mostafa, Dr. Madeha Saeed, Dr. Dalia Mohammad Abd Al-Salam, Dr. Nahed Mosaad"
$ws.Range("G90").Value = "Dr. Al-Shimaa Khaled, Dr. Mohammad Safwat"
$ws.Range("G95").Value = "Dr. Amr Saeed, Dr. Walaa Ghanima"
$ws.Range("G96").Value = "Dr. Ola Abd Al-Fattah, Dr. Eman Mohammad Al, Dr. Neveen Nashaat, Dr. Aya Hanafy, Dr. Salma Hassan, Dr. Marina Atef"
$ws.Range("G98").Value = "Dr. Mohammad El-Tanany, Dr. Rana Abo-Zaid, Dr. Servinaz Sayed Mohammad, Dr. Nesma, Dr. Nahla Nagiub"
$ws.Range("G100").Value = "Dr. Nourhan Mahmoud, Dr. Amira Sobhy, Dr. Alshimaa Atef, Dr. Eman Tantawi, Dr. Veronia Rafat"
$ws.Range("G101").Value = "Dr. Lamiaa Ossama, Dr. Nada Gouda, Dr. Fatma Elhady, Dr. Kerelos Zareef, Dr. Abeer Ragab, Dr. Nada Mohammad"
$ws.Range("G102").Value = "Dr. Fatma Elhady, Dr. Amera Ahmad Saad, Dr. Nada Gouda"
$ws.Range("G103").Value = "Dr. Arwa Elnagar, Dr. Shimaa Ashraf, Dr. Aya Saeed"
$ws.Range("G104").Value = "Dr. Eman M. Abo-Sakaya, Dr. Arwa Al-Sayed, Dr. Esraa Mostafa, Dr. Nourhan Osama, Dr. Yasmeena Fattoh, Dr. Marina Youhanna, Dr. Amira Ibrahim, Dr. Dina Adel, Dr. Maryam Ahmad"
$ws.Range("G105").Value = "Dr. Eman M. Abo-Sakaya, Dr. Esraa Mostafa, Dr. Nourhan Osama, Dr. Amany Raafat, Dr. Sarah Abdelmohsen, Dr. Merna Mahrous, Dr. Amira Ibrahim, Dr. Nadia Mostafa, Dr. Dina Adel, Dr. Maryam Ahmad"
$ws.Range("G112").Value = "Dr. Yassmen Ahmad, Dr. Nahla, Dr. Neveen Nashaat, Dr. Youstina Magdy, Dr. Salma Hassan, Dr. Remon"
$ws.Range("G114").Value = "Dr. Mohammad El-Tanany, Dr. Rana Abo-Zaid, Dr. Servinaz Sayed Mohammad, Dr. Nesma, Dr. Nahla Nagiub"
$ws.Range("G116").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Rana Abo-Zaid, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Gehan Adel, Dr. Majorelle Magdy, Dr. Manar Montaser"
$ws.Range("G117").Value = "Dr. Lamiaa Ossama, Dr. Nada Gouda, Dr. Fatma Elhady, Dr. Kerelos Zareef, Dr. Abeer Ragab, Dr. Nada Mohammad"
$ws.Range("G118").Value = "Dr. Fatma Elhady, Dr. Amera Ahmad Saad, Dr. Nada Gouda"
$ws.Range("G119").Value = "Dr. Arwa Elnagar, Dr. Shimaa Ashraf, Dr. Aya Saeed"
$ws.Range("G120").Value = "Dr. Eman M. Abo-Sakaya, Dr. Arwa Al-Sayed, Dr. Esraa Mostafa, Dr. Nourhan Osama, Dr. Yasmeena Fattoh, Dr. Marina Youhanna, Dr. Amira Ibrahim, Dr. Dina Adel, Dr. Maryam Ahmad"
$ws.Range("G121").Value = "Dr. Eman M. Abo-Sakaya, Dr. Esraa Mostafa, Dr. Nourhan Osama, Dr. Amany Raafat, Dr. Sarah Abdelmohsen, Dr. Merna Mahrous, Dr. Amira Ibrahim, Dr. Nadia Mostafa, Dr. Dina Adel, Dr. Maryam Ahmad"
$ws.Range("G125").Value = "Dr. Nancy Abd Al-Shafy, Dr. Walaa Ghanima"
$ws.Range("G128").Value = "Dr. Yassmen Ahmad, Dr. Nahla, Dr. Neveen Nashaat, Dr. Youstina Magdy, Dr. Salma Hassan, Dr. Remon"
